$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove a few stray cell values that were mistakenly populated by the
# naive component forecaster bug (C2, E2, C3 should be blank).
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()

# Corrected forecast values (tiny floating point fixes from the bug fix).
$ws.Range("C4").Value = 0.348613976222456
$ws.Range("C5").Value = -0.1384957661262676
$ws.Range("C6").Value = 1.566479473280191
$ws.Range("E6").Value = 2.49756057493542
$ws.Range("C7").Value = 0.7307568962937161
$ws.Range("C8").Value = 0.8188188121642126
$ws.Range("C10").Value = 1.9846842782967
$ws.Range("E10").Value = 2.047428048848809
$ws.Range("E11").Value = 1.552965246735782
$ws.Range("E12").Value = 0.232608152956959
$ws.Range("C13").Value = 1.064321453542272
$ws.Range("C14").Value = 1.361817904277718
$ws.Range("C15").Value = -4.352425014431327
$ws.Range("E15").Value = 31.54369540926345
$ws.Range("E16").Value = 22.41808675646531
$ws.Range("C18").Value = -0.9008525709169657
$ws.Range("E18").Value = -0.9756765446554017
$ws.Range("C19").Value = 0.2738544794132602
